$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M6").Value = 4392.44
$ws1.Range("M24").Value = "3 de 22"

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F6").Value = 4392.44
$ws2.Range("F24").Value = 8424.639999999999

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 7611.52
$ws3.Range("E16").Value = 28445.18
$ws3.Range("F16").Value = 0.2110986307676521
$ws3.Range("D19").Value = 8424.639999999999
$ws3.Range("E19").Value = 46598.52386304604
$ws3.Range("F19").Value = 0.1531107884121155
